$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New header for column D
$ws.Range("D1").Value = "PDRB"

# New PDRB data for rows 2-40 (one value per province row, plus the
# "Indonesia" total row 40)
$values = @(
    43782,
    73575,
    57047,
    165350,
    86722,
    75132,
    49233,
    51370,
    70194,
    161424,
    344350,
    56080,
    47972,
    51473,
    75770,
    70276,
    67319,
    32282,
    24272,
    52703,
    79320,
    67117,
    212175,
    198429,
    69352,
    120750,
    73573,
    67840,
    44433,
    42718,
    32198,
    70660,
    131636,
    59064,
    81009,
    61583,
    118774,
    18105,
    75020
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $values[$i]
}

# Restore the last recorded selection on Sheet1 (matches saved sheetView)
$ws.Range("F2").Select()
